$d = $word.ActiveDocument

# Remove the two chevron arrow shapes (Arrow: Chevron 320 and Arrow: Chevron 321).
# Note: in this environment Shapes.Item(index).Delete() operates on the shape at
# that position in underlying XML document order, and removing an item shifts
# later items down by one. The two chevrons are the 2nd and 3rd shapes in
# document order, so deleting position 2 twice removes both of them.
$d.Shapes.Item(2).Delete()
$d.Shapes.Item(2).Delete()
